# "Fruta / hortaliza, semanal" weekly data refresh.
#
# The data rows (2-16, excluding 6 and 11 which are unchanged) are
# re-shuffled: each row's Fecha/Volumen/Precio minimo/Precio maximo/
# Precio promedio ponderado/Precio $por Kg (columns D, J, K, L, M, P)
# are replaced with the values that another row held before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44400
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 9000
$ws.Range("L2").Value = 10000
$ws.Range("M2").Value = 9500
$ws.Range("P2").Value = 158
$ws.Range("D3").Value = 44421
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 9000
$ws.Range("M3").Value = 8500
$ws.Range("P3").Value = 142
$ws.Range("D4").Value = 44657
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 5000
$ws.Range("L4").Value = 5500
$ws.Range("M4").Value = 5250
$ws.Range("P4").Value = 88
$ws.Range("D5").Value = 44740
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 6000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 6500
$ws.Range("P5").Value = 108
$ws.Range("D7").Value = 44676
$ws.Range("J7").Value = 120
$ws.Range("K7").Value = 4000
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = 4250
$ws.Range("P7").Value = 71
$ws.Range("D8").Value = 44242
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 5000
$ws.Range("L8").Value = 5500
$ws.Range("M8").Value = 5250
$ws.Range("P8").Value = 88
$ws.Range("D9").Value = 44648
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 6750
$ws.Range("P9").Value = 112
$ws.Range("D10").Value = 44281
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 5500
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = 5750
$ws.Range("P10").Value = 96
$ws.Range("D12").Value = 44494
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 5000
$ws.Range("L12").Value = 6000
$ws.Range("M12").Value = 5500
$ws.Range("P12").Value = 92
$ws.Range("D13").Value = 44382
$ws.Range("J13").Value = 160
$ws.Range("K13").Value = 7000
$ws.Range("L13").Value = 8000
$ws.Range("M13").Value = 7438
$ws.Range("P13").Value = 124
$ws.Range("D14").Value = 44589
$ws.Range("J14").Value = 110
$ws.Range("K14").Value = 5000
$ws.Range("L14").Value = 6000
$ws.Range("M14").Value = 5500
$ws.Range("P14").Value = 92
$ws.Range("D15").Value = 44362
$ws.Range("J15").Value = 120
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 9000
$ws.Range("M15").Value = 8500
$ws.Range("P15").Value = 142
$ws.Range("D16").Value = 44627
$ws.Range("J16").Value = 120
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 4500
$ws.Range("M16").Value = 4250
$ws.Range("P16").Value = 71
